$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Update plan value in D5 (Residencial plan 50)
$ws.Range("D5").Value = 50

# Fill in row 6 (Residencial plan 100)
$ws.Range("B6").Value = "Residencial"
$ws.Range("C6").Value = "Sin_TotalPlay_TV"
$ws.Range("D6").Value = 100

# Fill in row 7 (Residencial plan 500)
$ws.Range("B7").Value = "Residencial"
$ws.Range("C7").Value = "Sin_TotalPlay_TV"
$ws.Range("D7").Value = 500

# Fill in row 8 (Residencial plan 1000)
$ws.Range("B8").Value = "Residencial"
$ws.Range("C8").Value = "Sin_TotalPlay_TV"
$ws.Range("D8").Value = 1000

# Reset the active selection on the sheet to A1
$ws.Activate()
$ws.Range("A1").Select()
